$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newText = "Profiles based on Bengaluru" + [char]10 + "See all profiles from" + [char]10 + "Bengaluru"

$ws.Range("D2").Value = $newText
$ws.Rows.Item(2).AutoFit()
